$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1434.4375
$ws.Range("I40").Value = 1307.8462
$ws.Range("J40").Value = 1521.0526
$ws.Range("K40").Value = 1307.8462
$ws.Range("L40").Value = 1521.0526
$ws.Range("M40").Value = -1132.8462
$ws.Range("N40").Value = -1871.0526

# Row 98
$ws.Range("H98").Value = 1965.75
$ws.Range("I98").Value = 1145.5
$ws.Range("K98").Value = 1145.5
$ws.Range("M98").Value = 352.5

# Row 122
$ws.Range("H122").Value = 1965.75
$ws.Range("I122").Value = 1145.5
$ws.Range("K122").Value = 3436.5
$ws.Range("M122").Value = -986.5

# Row 135
$ws.Range("H135").Value = 27752.816
$ws.Range("I135").Value = 34532.832
$ws.Range("K135").Value = 310795.488
$ws.Range("M135").Value = -308260.488

# Row 137
$ws.Range("H137").Value = 3659733
$ws.Range("I137").Value = 2001294.1
$ws.Range("J137").Value = 6251044
$ws.Range("K137").Value = 6003882.300000001
$ws.Range("L137").Value = 18753132
$ws.Range("M137").Value = -6001332.300000001
$ws.Range("N137").Value = -18758232

# Row 141
$ws.Range("H141").Value = 1526.6562
$ws.Range("I141").Value = 1024.909
$ws.Range("J141").Value = 2630.5
$ws.Range("K141").Value = 3074.727
$ws.Range("L141").Value = 7891.5
$ws.Range("M141").Value = 2105.273
$ws.Range("N141").Value = -18251.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17667.95
$ws.Range("I32").Value = 4460.9
$ws.Range("J32").Value = 79818.766
$ws.Range("K32").Value = 4460.9
$ws.Range("L32").Value = 79818.766
$ws.Range("M32").Value = -4173.9
$ws.Range("N32").Value = -80392.766

# Row 74
$ws.Range("H74").Value = 917.65
$ws.Range("I74").Value = 883.5294
$ws.Range("K74").Value = 883.5294
$ws.Range("M74").Value = -9.52940000000001

# Row 77
$ws.Range("H77").Value = 917.65
$ws.Range("I77").Value = 883.5294
$ws.Range("K77").Value = 4417.647
$ws.Range("M77").Value = -49.64699999999993

# Row 122
$ws.Range("H122").Value = 1538.8572
$ws.Range("I122").Value = 1531.3226
$ws.Range("J122").Value = 1560.091
$ws.Range("K122").Value = 4593.9678
$ws.Range("L122").Value = 4680.272999999999
$ws.Range("M122").Value = -2143.9678
$ws.Range("N122").Value = -9580.272999999999

# Row 132
$ws.Range("H132").Value = 127318.73
$ws.Range("I132").Value = 140370.81
$ws.Range("J132").Value = 9850
$ws.Range("K132").Value = 421112.43
$ws.Range("L132").Value = 29550
$ws.Range("M132").Value = -418582.43
$ws.Range("N132").Value = -34610

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2172.4348
$ws.Range("I86").Value = 2203.3
$ws.Range("K86").Value = 2203.3
$ws.Range("M86").Value = -1080.3

# Row 89
$ws.Range("H89").Value = 2172.4348
$ws.Range("I89").Value = 2203.3
$ws.Range("K89").Value = 11016.5
$ws.Range("M89").Value = -5400.5

# Row 134
$ws.Range("H134").Value = 63227.168
$ws.Range("I134").Value = 75603.83
$ws.Range("K134").Value = 226811.49
$ws.Range("M134").Value = -224276.49

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1329.62
$ws.Range("I58").Value = 1402.9111
$ws.Range("J58").Value = 670
$ws.Range("K58").Value = 1402.9111
$ws.Range("L58").Value = 670
$ws.Range("M58").Value = -1199.9111
$ws.Range("N58").Value = -1076

# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132
$ws.Range("H132").Value = 3274.5405
$ws.Range("I132").Value = 2823.8518
$ws.Range("K132").Value = 8471.555399999999
$ws.Range("M132").Value = -5941.555399999999

# Row 136
$ws.Range("H136").Value = 1329.62
$ws.Range("I136").Value = 1402.9111
$ws.Range("J136").Value = 670
$ws.Range("K136").Value = 4208.7333
$ws.Range("L136").Value = 2010
$ws.Range("M136").Value = -1658.7333
$ws.Range("N136").Value = -7110

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 945.25
$ws.Range("I131").Value = 493.33334
$ws.Range("J131").Value = 974.09576
$ws.Range("K131").Value = 1480.00002
$ws.Range("L131").Value = 2922.28728
$ws.Range("M131").Value = 3559.99998
$ws.Range("N131").Value = -13002.28728

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1708.2693
$ws.Range("I122").Value = 1532.5625
$ws.Range("J122").Value = 1989.4
$ws.Range("K122").Value = 4597.6875
$ws.Range("L122").Value = 5968.200000000001
$ws.Range("M122").Value = -2147.6875
$ws.Range("N122").Value = -10868.2

# Row 132
$ws.Range("H132").Value = 1689.1086
$ws.Range("I132").Value = 1299.9032
$ws.Range("J132").Value = 2493.4666
$ws.Range("K132").Value = 3899.7096
$ws.Range("L132").Value = 7480.399800000001
$ws.Range("M132").Value = -1369.7096
$ws.Range("N132").Value = -12540.3998

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 486.2143
$ws.Range("I16").Value = 477.25
$ws.Range("J16").Value = 498.16666
$ws.Range("K16").Value = 477.25
$ws.Range("L16").Value = 498.16666
$ws.Range("M16").Value = -307.25
$ws.Range("N16").Value = -838.16666

# Row 122
$ws.Range("H122").Value = 3201
$ws.Range("I122").Value = 2071.875
$ws.Range("J122").Value = 4491.4287
$ws.Range("K122").Value = 6215.625
$ws.Range("L122").Value = 13474.2861
$ws.Range("M122").Value = -3765.625
$ws.Range("N122").Value = -18374.2861

# Row 132
$ws.Range("H132").Value = 1501.6377
$ws.Range("I132").Value = 1424.3793
$ws.Range("K132").Value = 4273.1379
$ws.Range("M132").Value = -1743.1379

# Row 136
$ws.Range("H136").Value = 1102.3715
$ws.Range("I136").Value = 1017.4754
$ws.Range("J136").Value = 1677.7778
$ws.Range("K136").Value = 3052.4262
$ws.Range("L136").Value = 5033.3334
$ws.Range("M136").Value = -502.4261999999999
$ws.Range("N136").Value = -10133.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 6046002
$ws.Range("J11").Value = 3807502.5
$ws.Range("L11").Value = 3807502.5
$ws.Range("N11").Value = -3807786.5

# Row 81
$ws.Range("H81").Value = 4689.5186
$ws.Range("I81").Value = 5822.579
$ws.Range("J81").Value = 1998.5
$ws.Range("K81").Value = 11645.158
$ws.Range("L81").Value = 3997
$ws.Range("M81").Value = -10584.158
$ws.Range("N81").Value = -6119

# Row 84
$ws.Range("H84").Value = 4689.5186
$ws.Range("I84").Value = 5822.579
$ws.Range("J84").Value = 1998.5
$ws.Range("K84").Value = 58225.78999999999
$ws.Range("L84").Value = 19985
$ws.Range("M84").Value = -52921.78999999999
$ws.Range("N84").Value = -30593
